$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 8898.799999999999
$ws.Range("J125").Value = 3748.5
$ws.Range("L125").Value = 33736.5
$ws.Range("N125").Value = -38656.5
$ws.Range("H132").Value = 6104.1665
$ws.Range("I132").Value = 6060.9375
$ws.Range("K132").Value = 18182.8125
$ws.Range("M132").Value = -15652.8125
$ws.Range("H138").Value = 3654.75
$ws.Range("I138").Value = 1331.5294
$ws.Range("J138").Value = 4573.2324
$ws.Range("K138").Value = 3994.5882
$ws.Range("L138").Value = 13719.6972
$ws.Range("M138").Value = 1145.4118
$ws.Range("N138").Value = -23999.6972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2362.6528
$ws.Range("I32").Value = 2541.246
$ws.Range("J32").Value = 704.2857
$ws.Range("K32").Value = 2541.246
$ws.Range("L32").Value = 704.2857
$ws.Range("M32").Value = -2254.246
$ws.Range("N32").Value = -1278.2857
$ws.Range("H80").Value = 73999.5
$ws.Range("J80").Value = 73999.5
$ws.Range("L80").Value = 73999.5
$ws.Range("N80").Value = -75995.5
$ws.Range("H83").Value = 73999.5
$ws.Range("J83").Value = 73999.5
$ws.Range("L83").Value = 221998.5
$ws.Range("N83").Value = -231982.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14243.7
$ws.Range("I86").Value = 11745.333
$ws.Range("J86").Value = 17991.25
$ws.Range("K86").Value = 11745.333
$ws.Range("L86").Value = 17991.25
$ws.Range("M86").Value = -10622.333
$ws.Range("N86").Value = -20237.25
$ws.Range("H89").Value = 14243.7
$ws.Range("I89").Value = 11745.333
$ws.Range("J89").Value = 17991.25
$ws.Range("K89").Value = 58726.665
$ws.Range("L89").Value = 89956.25
$ws.Range("M89").Value = -53110.665
$ws.Range("N89").Value = -101188.25
$ws.Range("H96").Value = 27999
$ws.Range("I96").Value = 27999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 27999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -25253
$ws.Range("N96").ClearContents()
$ws.Range("H107").Value = 2620.926
$ws.Range("I107").Value = 2598.7144
$ws.Range("J107").Value = 2698.6667
$ws.Range("K107").Value = 2598.7144
$ws.Range("L107").Value = 2698.6667
$ws.Range("M107").Value = -678.7143999999998
$ws.Range("N107").Value = -6538.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1620.9
$ws.Range("I16").Value = 1675.8572
$ws.Range("J16").Value = 1492.6666
$ws.Range("K16").Value = 1675.8572
$ws.Range("L16").Value = 1492.6666
$ws.Range("M16").Value = -1388.8572
$ws.Range("N16").Value = -2066.6666
$ws.Range("H97").Value = 64000
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H107").Value = 15821.533
$ws.Range("I107").Value = 32080
$ws.Range("J107").Value = 1595.375
$ws.Range("K107").Value = 32080
$ws.Range("L107").Value = 1595.375
$ws.Range("M107").Value = -30160
$ws.Range("N107").Value = -5435.375
$ws.Range("H113").Value = 1620.9
$ws.Range("I113").Value = 1675.8572
$ws.Range("J113").Value = 1492.6666
$ws.Range("K113").Value = 1675.8572
$ws.Range("L113").Value = 1492.6666
$ws.Range("M113").Value = 494.1428000000001
$ws.Range("N113").Value = -5832.6666
$ws.Range("H132").Value = 15874.345
$ws.Range("I132").Value = 1450.1052
$ws.Range("K132").Value = 4350.3156
$ws.Range("M132").Value = -1820.3156
$ws.Range("H138").Value = 146427.58
$ws.Range("J138").Value = 146427.58
$ws.Range("L138").Value = 146427.58
$ws.Range("N138").Value = -156707.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 6207.923
$ws.Range("I81").Value = 1974.5
$ws.Range("J81").Value = 8089.4443
$ws.Range("K81").Value = 5923.5
$ws.Range("L81").Value = 24268.3329
$ws.Range("M81").Value = -4800.5
$ws.Range("N81").Value = -26514.3329
$ws.Range("H84").Value = 6207.923
$ws.Range("I84").Value = 1974.5
$ws.Range("J84").Value = 8089.4443
$ws.Range("K84").Value = 17770.5
$ws.Range("L84").Value = 72804.9987
$ws.Range("M84").Value = -12154.5
$ws.Range("N84").Value = -84036.9987
$ws.Range("H114").Value = 7636.3335
$ws.Range("J114").Value = 7636.3335
$ws.Range("L114").Value = 22909.0005
$ws.Range("N114").Value = -29417.0005
$ws.Range("H122").Value = 2483.138
$ws.Range("J122").Value = 2809.28
$ws.Range("L122").Value = 25283.52
$ws.Range("N122").Value = -30183.52
$ws.Range("H131").Value = 1458.4736
$ws.Range("I131").Value = 709.8
$ws.Range("K131").Value = 2129.4
$ws.Range("M131").Value = 2910.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10206.143
$ws.Range("I80").Value = 18668.334
$ws.Range("J80").Value = 3859.5
$ws.Range("K80").Value = 18668.334
$ws.Range("L80").Value = 3859.5
$ws.Range("M80").Value = -17670.334
$ws.Range("N80").Value = -5855.5
$ws.Range("H83").Value = 10206.143
$ws.Range("I83").Value = 18668.334
$ws.Range("J83").Value = 3859.5
$ws.Range("K83").Value = 93341.67
$ws.Range("L83").Value = 19297.5
$ws.Range("M83").Value = -88349.67
$ws.Range("N83").Value = -29281.5
$ws.Range("H126").Value = 20045.21
$ws.Range("J126").Value = 12999.363
$ws.Range("L126").Value = 38998.089
$ws.Range("N126").Value = -43938.089
$ws.Range("H132").Value = 2606.9092
$ws.Range("I132").Value = 2030.9231
$ws.Range("K132").Value = 6092.7693
$ws.Range("M132").Value = -3562.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 30609.75
$ws.Range("I22").Value = 57222
$ws.Range("J22").Value = 3997.5
$ws.Range("K22").Value = 57222
$ws.Range("L22").Value = 3997.5
$ws.Range("M22").Value = -56927
$ws.Range("N22").Value = -4587.5
$ws.Range("H24").Value = 15450.45
$ws.Range("I24").Value = 14875.625
$ws.Range("K24").Value = 14875.625
$ws.Range("M24").Value = -14532.625
$ws.Range("H27").Value = 30609.75
$ws.Range("I27").Value = 57222
$ws.Range("J27").Value = 3997.5
$ws.Range("K27").Value = 57222
$ws.Range("L27").Value = 3997.5
$ws.Range("M27").Value = -57115
$ws.Range("N27").Value = -4211.5
$ws.Range("H74").Value = 80000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 80000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H132").Value = 483497.3
$ws.Range("I132").Value = 785363.6
$ws.Range("K132").Value = 2356090.8
$ws.Range("M132").Value = -2353560.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 34996.668
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H62").Value = 213596.06
$ws.Range("I62").Value = 544532.7
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 544532.7
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -543908.7
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 213596.06
$ws.Range("I65").Value = 544532.7
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 2722663.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -2719543.5
$ws.Range("N65").Value = -21240
$ws.Range("H100").Value = 32215.475
$ws.Range("I100").Value = 19472.934
$ws.Range("K100").Value = 38945.868
$ws.Range("M100").Value = -38404.868
